# TFS 12308 - New sub coaching reason for warnings
# 1) Revision_History: log the change
# 2) DIM_Sub_Coaching_Reason: add the 3 new sub-coaching reasons (highlight new, un-highlight old)
# 3) Coaching_Reason_Selection: add the 9 new reason/sub-reason combinations (highlight new, un-highlight old)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Revision_History (sheet "Revision_History") - append row 73
# ---------------------------------------------------------------------------
$wsHist = $wb.Worksheets.Item("Revision_History")

$wsHist.Range("A73").Value = 64
$wsHist.Range("B72").Copy()
$wsHist.Range("B73").PasteSpecial(-4122)
$wsHist.Range("B73").Value = 43383
$wsHist.Range("C73").Value = "Susmitha Palacherla"
$wsHist.Range("D73").Value = 12308
$wsHist.Range("E73").Value = "New sub coaching reason for warnings. Added records to SubCoaching Reason and Coaching Reason Selection table tabs"
[void]$wsHist.Range("E73").Select()

# ---------------------------------------------------------------------------
# 2) DIM_Sub_Coaching_Reason (sheet "DIM_Sub_Coaching_Reason") - append rows 248:250
# ---------------------------------------------------------------------------
$wsSub = $wb.Worksheets.Item("DIM_Sub_Coaching_Reason")

# Copy the highlighted format of the previous "newest" row down onto the new rows,
# then clear the highlight off the old row (matches the diff's style move).
$wsSub.Range("A247:B247").Copy()
$wsSub.Range("A248:B250").PasteSpecial(-4122)
$wsSub.Range("A247:B247").ClearFormats()

$subRows = @(
    @(246, "Adherence"),
    @(247, "Quality/Performance – Failed Calls"),
    @(248, "Quality/Performance – Critical Fails")
)

$r = 248
foreach ($row in $subRows) {
    $wsSub.Cells.Item($r, 1).Value = $row[0]
    $wsSub.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

[void]$wsSub.Range("B253").Select()

# ---------------------------------------------------------------------------
# 3) Coaching_Reason_Selection (sheet "Coaching_Reason_Selection") - append rows 310:318
# ---------------------------------------------------------------------------
$wsSel = $wb.Worksheets.Item("Coaching_Reason_Selection")

# Copy the highlighted format of the previous "newest" row down onto the new rows,
# then clear the highlight off the old row (matches the diff's style move).
$wsSel.Range("A309:P309").Copy()
$wsSel.Range("A310:P318").PasteSpecial(-4122)
$wsSel.Range("A309:P309").ClearFormats()

$flags = @(1,1,0,1,0,1,1,1,1,1,1,1)

$selRows = @(
    @(28, "Verbal Warning", 246, "Adherence"),
    @(28, "Verbal Warning", 247, "Quality/Performance – Failed Calls"),
    @(28, "Verbal Warning", 248, "Quality/Performance – Critical Fails"),
    @(29, "Written Warning", 246, "Adherence"),
    @(29, "Written Warning", 247, "Quality/Performance – Failed Calls"),
    @(29, "Written Warning", 248, "Quality/Performance – Critical Fails"),
    @(30, "Final Written Warning", 246, "Adherence"),
    @(30, "Final Written Warning", 247, "Quality/Performance – Failed Calls"),
    @(30, "Final Written Warning", 248, "Quality/Performance – Critical Fails")
)

$r = 310
foreach ($row in $selRows) {
    $wsSel.Cells.Item($r, 1).Value = $row[0]
    $wsSel.Cells.Item($r, 2).Value = $row[1]
    $wsSel.Cells.Item($r, 3).Value = $row[2]
    $wsSel.Cells.Item($r, 4).Value = $row[3]
    for ($c = 0; $c -lt $flags.Length; $c++) {
        $wsSel.Cells.Item($r, 5 + $c).Value = $flags[$c]
    }
    $r = $r + 1
}

[void]$wsSel.Range("A319").Select()
